$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("cases")
$sheet2 = $wb.Worksheets.Item("deaths")

$sheet1.Range("P15").Value = 35986
$sheet1.Range("P16").Value = 38806
$sheet1.Range("P17").Value = 40430
$sheet1.Range("P18").Value = 42796
$sheet1.Range("P19").Value = 45553
$sheet1.Range("P20").Value = 49995
$sheet1.Range("P21").Value = 52238
$sheet1.Range("P22").Value = 54328
$sheet1.Range("P23").Value = 56728
$sheet1.Range("P24").Value = 60698
$sheet1.Range("P25").Value = 63272
$sheet1.Range("P26").Value = 66716
$sheet1.Range("P27").Value = 72022
$sheet1.Range("P28").Value = 75128

$sheet2.Range("P15").Value = 2363
$sheet2.Range("P16").Value = 2605
$sheet2.Range("P17").Value = 2747
$sheet2.Range("P18").Value = 2956
$sheet2.Range("P19").Value = 3205
$sheet2.Range("P20").Value = 3614
$sheet2.Range("P21").Value = 3825
$sheet2.Range("P22").Value = 4023
$sheet2.Range("P23").Value = 4254
$sheet2.Range("P24").Value = 4643
$sheet2.Range("P25").Value = 4899
$sheet2.Range("P26").Value = 5246
$sheet2.Range("P27").Value = 5791
$sheet2.Range("P28").Value = 6115
